$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Apply updated crypto market data values. Cells whose new value would
# otherwise be auto-parsed by Excel as a number are temporarily forced to
# Text format so they are stored as strings (matching the source data),
# then restored to the default "Normal" style so no style override remains.
$ws.Range("D2").Value = "26.955.98"
$ws.Range("E2").Value = "  +1.58%  "
$ws.Range("D3").Value = "1.672.74"
$ws.Range("E3").Value = "  +2.93%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("E6").Value = "  +6.21%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("E8").Value = "  +2.71%  "
$ws.Range("E9").Value = "  +1.86%  "
$ws.Range("E10").Value = "  +5.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0891"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.07%  "
$ws.Range("D12").Value = "1.909.92"
$ws.Range("E12").Value = "  +3.11%  "
$ws.Range("D13").Value = "1.672.16"
$ws.Range("E13").Value = "  +3.11%  "
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.522"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.85%  "
$ws.Range("D17").Value = "26.979.61"
$ws.Range("E17").Value = "  +1.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "232.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("E19").Value = "  +1.28%  "
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("E22").Value = "  +2.93%  "
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.117"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.94%  "
$ws.Range("E27").Value = "  +0.94%  "
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("E31").Value = "  +1.46%  "
$ws.Range("E32").Value = "  +1.79%  "
$ws.Range("D33").Value = "1.455.10"
$ws.Range("E33").Value = "  -4.64%  "
$ws.Range("E34").Value = "  +4.44%  "
$ws.Range("E35").Value = "  +4.40%  "
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("E37").Value = "  +7.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.566"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("E39").Value = "  +1.38%  "
$ws.Range("E40").Value = "  +3.34%  "
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("E42").Value = "  +4.06%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.81%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.964"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.36%  "
$ws.Range("D45").Value = "1.815.55"
$ws.Range("E45").Value = "  +2.88%  "
$ws.Range("E46").Value = "  +2.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.88%  "
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("E49").Value = "  +0.78%  "
$ws.Range("E50").Value = "  +4.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0508"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.21%  "
